$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing header labels right by one column to make room for the new
# "Quantity:" column, then write the new label into C1.
$ws.Range("F1").Value = $ws.Range("E1").Value()
$ws.Range("E1").Value = $ws.Range("D1").Value()
$ws.Range("D1").Value = $ws.Range("C1").Value()
$ws.Range("C1").Value = "Quantity:"

# Set the width of the new column F (previously used by Total Cost:)
# (ColumnWidth input gets a +5/6 padding offset applied when stored, so back
# it off here to land on a stored width of exactly 10.)
$ws.Columns.Item(6).ColumnWidth = 9.166666666666666

# Move the selection to C2, matching the saved cursor position
$ws.Range("C2").Select()
